$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "28.524.71"
$ws.Range("D3").Value = "1.825.43"
Set-TextValue $ws "D5" "317.36"
Set-TextValue $ws "D7" "0.5188"
Set-TextValue $ws "D9" "0.08443"
Set-TextValue $ws "D12" "6.426"
Set-TextValue $ws "D13" "21.09"
Set-TextValue $ws "D15" "7.535"
$ws.Range("D16").Value = "1.822.19"
Set-TextValue $ws "D17" "0.00001134"
Set-TextValue $ws "D18" "93.03"
Set-TextValue $ws "D19" "0.06599"
Set-TextValue $ws "D20" "17.79"
Set-TextValue $ws "D22" "6.073"
$ws.Range("D23").Value = "28.546.18"
Set-TextValue $ws "D24" "11.44"
Set-TextValue $ws "D26" "21.09"
Set-TextValue $ws "D27" "159.60"
$ws.Range("D28").Value = "2.030.59"
Set-TextValue $ws "D29" "2.404"
Set-TextValue $ws "D30" "125.67"
Set-TextValue $ws "D32" "1.101"
Set-TextValue $ws "D33" "5.729"
Set-TextValue $ws "D34" "0.07439"
Set-TextValue $ws "D35" "3.660"
Set-TextValue $ws "D36" "0.2229"
Set-TextValue $ws "D37" "0.02361"
Set-TextValue $ws "D38" "5.232"
Set-TextValue $ws "D39" "8.821"
Set-TextValue $ws "D40" "0.6329"
Set-TextValue $ws "D41" "11.38"
Set-TextValue $ws "D42" "1.193"
Set-TextValue $ws "D43" "1.401"
Set-TextValue $ws "D44" "13.60"
Set-TextValue $ws "D45" "3.783"
Set-TextValue $ws "D46" "0.5973"
Set-TextValue $ws "D47" "126.65"
Set-TextValue $ws "D50" "0.06986"
Set-TextValue $ws "D51" "74.46"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("E9").Value = "  +8.85%  "
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  +4.33%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  -0.10%  "

# --- Row 10 / Row 11 swap (OKB <-> Polygon) ---
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws "D10" "1.115"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D11" "41.89"
$ws.Range("E11").Value = "  -0.32%  "
